$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.193.51'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '1.859.82'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7145'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07749'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3076'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08256'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.231'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.847.47'
$ws.Range("E13").Value = '  -1.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7141'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '29.204.51'
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.867'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '244.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.112.66'
$ws.Range("E21").Value = '  -1.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.937'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1583'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.918'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.493'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.314'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.385'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.142'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.907'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.174'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7278'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.678'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01848'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.687'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").Value = '1.154.70'
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9047'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.096'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").Value = '2.007.88'
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5236'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.765'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.304'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.868'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.69%  '
